$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(130).Insert()

$ws.Range("A130").Value = 7
$ws.Range("B130").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C130").Value = "Ñuble"
$ws.Range("D130").Value = 45223
$ws.Range("E130").Value = 16
$ws.Range("F130").Value = 100112013
$ws.Range("G130").Value = "Alcachofa"
$ws.Range("H130").Value = "Argentina(o)"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 60
$ws.Range("K130").Value = 13000
$ws.Range("L130").Value = 13000
$ws.Range("M130").Value = 13000
$ws.Range("N130").Value = "$/caja 50 unidades"
$ws.Range("O130").Value = "Provincia de Limarí"
$ws.Range("P130").Value = 260
$ws.Range("Q130").Value = 50
$ws.Range("R130").Value = "Hortaliza"

Write-Host $ws.UsedRange.Address()
